# Remove the Jekyll/Github-pages copyright footer block, along with the
# blank spacer paragraph and the page-break paragraph that introduced it,
# which previously followed the "LOQ4056: ..." requirements line.
#
# Paragraph layout before the edit (end of document):
#   ... "LOQ4056: Química Analítica para Engenharia (Requisito fraco)"
#   <empty paragraph>                                   <- delete
#   <empty paragraph, page-break-before, jc=left>        <- delete
#   "© 2020 . Contact: luizeleno@usp.br. ..."            <- delete
#   <empty paragraph>                                    <- keep
#   <empty paragraph, page-break-before>                 <- keep

$d = $word.ActiveDocument

$anchor = "LOQ4056: Química Analítica para Engenharia (Requisito fraco)"

# Find the paragraph containing the anchor text and walk forward from the
# paragraph right after it.
$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($paragraphs.Item($i).Range.Text -like "*$anchor*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    # The copyright text identifies the last paragraph of the block we
    # want removed; find it among the paragraphs following the anchor.
    $copyrightIndex = -1
    for ($j = $anchorIndex + 1; $j -le $count; $j++) {
        if ($paragraphs.Item($j).Range.Text -like "*Powered by Jekyll*") {
            $copyrightIndex = $j
            break
        }
    }

    if ($copyrightIndex -gt 0) {
        $startRange = $paragraphs.Item($anchorIndex + 1).Range.Start
        $endRange = $paragraphs.Item($copyrightIndex).Range.End
        $d.Range($startRange, $endRange).Delete()
    }
}
